$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 325, pushing existing rows 325-392 down to 327-394.
$ws.Rows("325:326").Insert()

# --- New row 325 ---
$ws.Cells.Item(325, 1).Value = 5
$ws.Cells.Item(325, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(325, 3).Value = "Maule"
$ws.Cells.Item(325, 4).Value = 44798
$ws.Cells.Item(325, 5).Value = 7
$ws.Cells.Item(325, 6).Value = 100112023
$ws.Cells.Item(325, 7).Value = "Brócoli"
$ws.Cells.Item(325, 8).Value = "Sin especificar"
$ws.Cells.Item(325, 9).Value = "Primera"
$ws.Cells.Item(325, 10).Value = 2000
$ws.Cells.Item(325, 11).Value = 900
$ws.Cells.Item(325, 12).Value = 900
$ws.Cells.Item(325, 13).Value = 900
$ws.Cells.Item(325, 14).Value = "`$/unidad"
$ws.Cells.Item(325, 15).Value = "Región del Maule"
$ws.Cells.Item(325, 16).Value = 900
$ws.Cells.Item(325, 17).Value = 1
$ws.Cells.Item(325, 18).Value = "Hortaliza"

# --- New row 326 ---
$ws.Cells.Item(326, 1).Value = 5
$ws.Cells.Item(326, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(326, 3).Value = "Maule"
$ws.Cells.Item(326, 4).Value = 44798
$ws.Cells.Item(326, 5).Value = 7
$ws.Cells.Item(326, 6).Value = 100112023
$ws.Cells.Item(326, 7).Value = "Brócoli"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Segunda"
$ws.Cells.Item(326, 10).Value = 3000
$ws.Cells.Item(326, 11).Value = 700
$ws.Cells.Item(326, 12).Value = 700
$ws.Cells.Item(326, 13).Value = 700
$ws.Cells.Item(326, 14).Value = "`$/unidad"
$ws.Cells.Item(326, 15).Value = "Región del Maule"
$ws.Cells.Item(326, 16).Value = 700
$ws.Cells.Item(326, 17).Value = 1
$ws.Cells.Item(326, 18).Value = "Hortaliza"
